# Updates cryptos list values (Price and Volume(1h) columns) to match the
# latest scrape, per commit "Updated cryptos list on Sat Oct 19 13:30:50 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds free-form numeric-looking text (e.g. "68.156.48",
# "2.638.22", "1.00", "0.0000192"). Force Text format so Excel keeps the exact
# characters (no auto-conversion to a Double / scientific notation / dropped zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.156.48'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.638.22'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.01'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.50'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.545'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.637.65'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.351'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000192'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.84'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.120.01'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.108.27'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.637.33'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '362.79'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.81'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.73'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.68'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.773.25'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '561.45'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.98'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.85'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.129'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.34'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.30'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.372'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₆0340'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.46'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '158.13'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.93'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0785'

# Column E ("Volume(1h)") values are padded percentage strings; plain text assignment
# keeps them as-is (the leading "+"/"-" with "%" prevents Excel from treating them as numbers).
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("E9").Value = '  +0.70%  '
$ws.Range("E10").Value = '  +8.64%  '
$ws.Range("E11").Value = '  -0.88%  '
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("E13").Value = '  +1.41%  '
$ws.Range("E14").Value = '  +3.11%  '
$ws.Range("E15").Value = '  +1.02%  '
$ws.Range("E16").Value = '  +0.78%  '
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("E20").Value = '  -1.32%  '
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("E22").Value = '  +3.51%  '
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("E25").Value = '  +2.83%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("E27").Value = '  -2.53%  '
$ws.Range("E28").Value = '  +1.89%  '
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("E31").Value = '  -1.26%  '
$ws.Range("E32").Value = '  +1.45%  '
$ws.Range("E33").Value = '  +0.57%  '
$ws.Range("E34").Value = '  +1.34%  '
$ws.Range("E35").Value = '  +2.11%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("E37").Value = '  +3.65%  '
$ws.Range("E38").Value = '  -0.88%  '
$ws.Range("E39").Value = '  +1.15%  '
$ws.Range("E40").Value = '  +1.53%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("E42").Value = '  +0.33%  '
$ws.Range("E43").Value = '  +3.80%  '
$ws.Range("E44").Value = '  +2.52%  '
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("E47").Value = '  +0.76%  '
$ws.Range("E48").Value = '  +2.14%  '
$ws.Range("E49").Value = '  +1.89%  '
$ws.Range("E50").Value = '  +1.02%  '
$ws.Range("E51").Value = '  +1.42%  '
